# Update the 100 arithmetic problems in the table (20 rows x 5 columns)
# to match the new values, in row-major (left-to-right, top-to-bottom) order.

$newValues = @(
    "39+33=",
    "7+38=",
    "57+6=",
    "32-16=",
    "27+64=",
    "64-7=",
    "74+18=",
    "78+4=",
    "70-44=",
    "49+24=",
    "37+44=",
    "4+57=",
    "94-25=",
    "7+88=",
    "60-33=",
    "80-79=",
    "47+37=",
    "65-6=",
    "82-55=",
    "39+22=",
    "29+49=",
    "9+85=",
    "65+16=",
    "48+26=",
    "74-8=",
    "61-55=",
    "19+15=",
    "63-5=",
    "74-36=",
    "3+78=",
    "38+37=",
    "86+8=",
    "75+17=",
    "73-39=",
    "53-15=",
    "93-5=",
    "79+19=",
    "75+16=",
    "19+62=",
    "75-68=",
    "8+33=",
    "97-78=",
    "9+12=",
    "9+38=",
    "59+13=",
    "57+7=",
    "39+29=",
    "73-46=",
    "51-9=",
    "86-57=",
    "18+74=",
    "58+34=",
    "29+22=",
    "59+8=",
    "45+17=",
    "81-72=",
    "27+28=",
    "91-7=",
    "9+57=",
    "73-27=",
    "62-27=",
    "16+55=",
    "92-63=",
    "45-38=",
    "9+32=",
    "47+25=",
    "16+5=",
    "34+38=",
    "38+4=",
    "35-8=",
    "84-69=",
    "56-37=",
    "37+16=",
    "79+19=",
    "47+6=",
    "58+9=",
    "94-45=",
    "37+36=",
    "60-54=",
    "72-63=",
    "72-5=",
    "13-9=",
    "18+25=",
    "16+68=",
    "19+14=",
    "90-44=",
    "34+18=",
    "52-34=",
    "85-69=",
    "61-15=",
    "87+7=",
    "15+6=",
    "37+36=",
    "7+54=",
    "36+29=",
    "32-6=",
    "51-33=",
    "48+36=",
    "93-44=",
    "97-28="
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.MoveEnd(1, -1) | Out-Null
        $rng.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
